$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the changed rows.
# D-column text that looks like a plain number is written via a Text
# number-format round-trip so Excel keeps it as the literal string
# (e.g. "212.20", "0.572") instead of silently parsing it into a float
# and losing the trailing zero / exact decimal text. The format is reset
# back to Normal afterwards so no stray style sticks to the cell.

$ws.Range("D2").Value = "27.458.52"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "1.636.58"
$ws.Range("E3").Value = "  -0.68%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("E6").Value = "  +4.68%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.55%  "

$ws.Range("E9").Value = "  -2.24%  "

$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").Value = "1.869.21"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").Value = "1.636.69"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.572"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.00%  "

$ws.Range("E15").Value = "  -1.67%  "

$ws.Range("E16").Value = "  -2.37%  "

$ws.Range("D17").Value = "27.474.38"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.60%  "

$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("E20").Value = "  +2.88%  "

$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.77%  "

$ws.Range("E24").Value = "  -3.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.112"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.53%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("E29").Value = "  -3.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.52%  "

$ws.Range("E31").Value = "  -1.73%  "

$ws.Range("E32").Value = "  -0.39%  "

$ws.Range("E33").Value = "  +2.94%  "

$ws.Range("D34").Value = "1.412.27"
$ws.Range("E34").Value = "  -3.02%  "

$ws.Range("E35").Value = "  +2.70%  "

$ws.Range("E36").Value = "  -2.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.570"
$ws.Range("D37").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.923"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +18.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.874"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.68%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.72%  "

$ws.Range("E44").Value = "  +1.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.63%  "

$ws.Range("D46").Value = "1.777.68"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("E47").Value = "  -3.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.66%  "

$ws.Range("E49").Value = "  +0.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0987"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.56%  "
